# Update the Event Data Flow Table:
# - Clear the content of the "query" rows (13 and 14), leaving the
#   cell formatting/borders in place (style s="1" is preserved).
# - Move the active selection to F17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents (not formatting) for rows 13 and 14, columns A:F
$ws.Range("A13:F14").ClearContents()

# Update the active cell selection to F17
$ws.Range("F17").Select()
